# Commit: "Wed, May 27, 2020  8:05:21 AM"
#
# The deck currently has the "Integral" design applied (theme colours on
# ppt/theme/theme2.xml, the theme wired to the SlideMaster / Presentation).
# The edit reverts the applied design's colour scheme back to the stock
# "Office Theme" palette (the palette that this deck's Notes Master /
# theme1.xml already carries), i.e. the twelve theme colour slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) of the presentation's
# live theme are switched from the Integral colours to the Office colours.
#
# PowerPoint's ColorScheme/RGBColor COM objects store RGB as a single
# Long in B + G*256 + R*65536 order, so a literal 0xBBGGRR below reads as
# the familiar RRGGBB hex value.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$scheme = $theme.ThemeColorScheme

# 1 dk1
$scheme.Item(1).RGB = 0x000000
# 2 lt1
$scheme.Item(2).RGB = 0xFFFFFF
# 3 dk2
$scheme.Item(3).RGB = 0x6A5444
# 4 lt2
$scheme.Item(4).RGB = 0xE6E6E7
# 5 accent1
$scheme.Item(5).RGB = 0xD59B5B
# 6 accent2
$scheme.Item(6).RGB = 0x317DED
# 7 accent3
$scheme.Item(7).RGB = 0xA5A5A5
# 8 accent4
$scheme.Item(8).RGB = 0x00C0FF
# 9 accent5
$scheme.Item(9).RGB = 0xC47244
# 10 accent6
$scheme.Item(10).RGB = 0x47AD70
# 11 hlink
$scheme.Item(11).RGB = 0xC16305
# 12 folHlink
$scheme.Item(12).RGB = 0x724F95
